$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.276.66"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.472.55"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.83"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.472.07"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.69"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.072.04"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.49"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.84%  "
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.479.50"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.364.34"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("E19").Value = "  +3.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.33"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.76"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.621.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("E28").Value = "  -4.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.62"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.24"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.53"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.34"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.16"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.71%  "
$ws.Range("E38").Value = "  +8.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.79"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +11.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "169.10"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.515.28"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0763"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.799"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.73"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.37"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  +2.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.600.99"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.31"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +9.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.09"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +2.72%  "
